# Applies the cryptos.xlsx price/volume/row-reorder update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-three glyph (U+2083) used in the PEPE unit price; built from the char code
# so the source file itself can stay plain ASCII.
$sub3 = [string][char]8323

$ws.Range("D2").Value = '61.198.51'
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = '2.390.35'
$ws.Range("E3").Value = '  -4.08%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''547.89'
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("D6").Value = '''142.24'
$ws.Range("E6").Value = '  -3.74%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '''0.539'
$ws.Range("E8").Value = '  -10.95%  '
$ws.Range("D9").Value = '2.388.56'
$ws.Range("E9").Value = '  -4.09%  '
$ws.Range("E10").Value = '  -2.37%  '
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '''5.25'
$ws.Range("E12").Value = '  -3.89%  '
$ws.Range("D13").Value = '''0.348'
$ws.Range("E13").Value = '  -3.60%  '
$ws.Range("D14").Value = '''25.40'
$ws.Range("E14").Value = '  -3.52%  '
$ws.Range("D15").Value = '2.818.47'
$ws.Range("E15").Value = '  -4.19%  '
$ws.Range("E16").Value = '  -2.75%  '
$ws.Range("D17").Value = '61.077.31'
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").Value = '2.387.89'
$ws.Range("E18").Value = '  -4.36%  '
$ws.Range("E19").Value = '  -4.18%  '
$ws.Range("E20").Value = '  -1.91%  '
$ws.Range("D21").Value = '''318.91'
$ws.Range("E21").Value = '  -1.56%  '
$ws.Range("D22").Value = '''6.75'
$ws.Range("E22").Value = '  -4.26%  '
$ws.Range("E23").Value = '  +8.32%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = '''63.81'
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("D26").Value = '''8.17'
$ws.Range("E26").Value = '  +7.14%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("B28").Value = 'Bittensor'
$ws.Range("C28").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D28").Value = '''538.58'
$ws.Range("E28").Value = '  -0.91%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0' + $sub3 + '0940'
$ws.Range("E29").Value = '  -6.96%  '
$ws.Range("D30").Value = '2.506.19'
$ws.Range("E30").Value = '  -4.00%  '
$ws.Range("E31").Value = '  -5.97%  '
$ws.Range("D32").Value = '''8.12'
$ws.Range("E32").Value = '  -3.82%  '
$ws.Range("E33").Value = '  -4.12%  '
$ws.Range("E34").Value = '  -3.55%  '
$ws.Range("D35").Value = '''1.59'
$ws.Range("E35").Value = '  -1.71%  '
$ws.Range("D36").Value = '''0.998'
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").Value = '''5.59'
$ws.Range("E37").Value = '  -7.21%  '
$ws.Range("D38").Value = '''4.74'
$ws.Range("E38").Value = '  -4.01%  '
$ws.Range("D39").Value = '''0.378'
$ws.Range("E39").Value = '  -2.48%  '
$ws.Range("D40").Value = '''1.83'
$ws.Range("E40").Value = '  +4.20%  '
$ws.Range("D41").Value = '''18.13'
$ws.Range("E41").Value = '  -2.67%  '
$ws.Range("D42").Value = '''138.93'
$ws.Range("E42").Value = '  -6.55%  '
$ws.Range("D44").Value = '''40.34'
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").Value = '''2.23'
$ws.Range("E45").Value = '  -6.91%  '
$ws.Range("D46").Value = '''142.08'
$ws.Range("E46").Value = '  -4.98%  '
$ws.Range("E47").Value = '  -0.81%  '
$ws.Range("D48").Value = '''20.27'
$ws.Range("E48").Value = '  -3.89%  '
$ws.Range("D49").Value = '''0.0520'
$ws.Range("E49").Value = '  -3.12%  '
$ws.Range("D50").Value = '''0.578'
$ws.Range("E50").Value = '  -3.66%  '
$ws.Range("D51").Value = '''0.0227'
$ws.Range("E51").Value = '  -1.44%  '
